# Apply cell value updates as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 10).Value = 1.03  # J4: 1.04 -> 1.03
$ws.Cells.Item(4, 11).Value = 15  # K4: 13 -> 15
$ws.Cells.Item(4, 12).Value = 1.2  # L4: 1.22 -> 1.2
$ws.Cells.Item(4, 13).Value = 4.33  # M4: 4 -> 4.33
# Row 6
$ws.Cells.Item(6, 9).Value = 6.6  # I6: 6.5 -> 6.6
$ws.Cells.Item(6, 25).Value = 32  # Y6: 35 -> 32
$ws.Cells.Item(6, 32).Value = 40  # AF6: 37 -> 40
$ws.Cells.Item(6, 34).Value = 175  # AH6: 150 -> 175
# Row 7
$ws.Cells.Item(7, 20).Value = 6.3  # T7: 6.2 -> 6.3
$ws.Cells.Item(7, 25).Value = 29  # Y7: 30 -> 29
$ws.Cells.Item(7, 27).Value = 6.9  # AA7: 7 -> 6.9
# Row 9
$ws.Cells.Item(9, 8).Value = 3.4  # H9: 3.35 -> 3.4
$ws.Cells.Item(9, 9).Value = 4.05  # I9: 4.2 -> 4.05
$ws.Cells.Item(9, 12).Value = 1.33  # L9: 1.37 -> 1.33
$ws.Cells.Item(9, 13).Value = 2.8  # M9: 2.65 -> 2.8
$ws.Cells.Item(9, 14).Value = 1.98  # N9: 2.07 -> 1.98
$ws.Cells.Item(9, 15).Value = 1.65  # O9: 1.6 -> 1.65
$ws.Cells.Item(9, 16).Value = 1.42  # P9: 1.47 -> 1.42
$ws.Cells.Item(9, 17).Value = 2.45  # Q9: 2.35 -> 2.45
$ws.Cells.Item(9, 18).Value = 1.87  # R9: 1.93 -> 1.87
$ws.Cells.Item(9, 19).Value = 1.75  # S9: 1.7 -> 1.75
$ws.Cells.Item(9, 20).Value = 6.3  # T9: 5.9 -> 6.3
$ws.Cells.Item(9, 21).Value = 8  # U9: 7.7 -> 8
$ws.Cells.Item(9, 24).Value = 15.5  # X9: 16 -> 15.5
$ws.Cells.Item(9, 25).Value = 30  # Y9: 32 -> 30
$ws.Cells.Item(9, 26).Value = 9  # Z9: 8.25 -> 9
$ws.Cells.Item(9, 27).Value = 6.7  # AA9: 6.5 -> 6.7
$ws.Cells.Item(9, 28).Value = 16.5  # AB9: 17.5 -> 16.5
$ws.Cells.Item(9, 29).Value = 90  # AC9: 100 -> 90
$ws.Cells.Item(9, 30).Value = 800  # AD9: 1000 -> 800
$ws.Cells.Item(9, 31).Value = 10.5  # AE9: 10 -> 10.5
$ws.Cells.Item(9, 33).Value = 14  # AG9: 14.5 -> 14
$ws.Cells.Item(9, 34).Value = 65  # AH9: 70 -> 65
$ws.Cells.Item(9, 35).Value = 40  # AI9: 45 -> 40
$ws.Cells.Item(9, 36).Value = 50  # AJ9: 55 -> 50
# Row 12
$ws.Cells.Item(12, 7).Value = 3.5  # G12: 3.6 -> 3.5
$ws.Cells.Item(12, 8).Value = 2.42  # H12: 2.45 -> 2.42
$ws.Cells.Item(12, 9).Value = 2.62  # I12: 2.55 -> 2.62
$ws.Cells.Item(12, 11).Value = 4.35  # K12: 4.4 -> 4.35
$ws.Cells.Item(12, 13).Value = 2.1  # M12: 2.12 -> 2.1
$ws.Cells.Item(12, 14).Value = 2.9  # N12: 2.87 -> 2.9
$ws.Cells.Item(12, 15).Value = 1.36  # O12: 1.37 -> 1.36
$ws.Cells.Item(12, 17).Value = 2.1  # Q12: 2.12 -> 2.1
$ws.Cells.Item(12, 20).Value = 7  # T12: 7.2 -> 7
$ws.Cells.Item(12, 21).Value = 17  # U12: 18 -> 17
$ws.Cells.Item(12, 22).Value = 12.5  # V12: 13 -> 12.5
$ws.Cells.Item(12, 26).Value = 4.35  # Z12: 4.4 -> 4.35
$ws.Cells.Item(12, 28).Value = 17.5  # AB12: 17 -> 17.5
$ws.Cells.Item(12, 31).Value = 5.7  # AE12: 5.6 -> 5.7
$ws.Cells.Item(12, 32).Value = 11.5  # AF12: 11 -> 11.5
$ws.Cells.Item(12, 34).Value = 32  # AH12: 30 -> 32
$ws.Cells.Item(12, 35).Value = 30  # AI12: 28 -> 30
# Row 14
$ws.Cells.Item(14, 7).Value = 2.42  # G14: 2.3 -> 2.42
$ws.Cells.Item(14, 9).Value = 3  # I14: 3.2 -> 3
$ws.Cells.Item(14, 17).Value = 2.37  # Q14: 2.35 -> 2.37
$ws.Cells.Item(14, 20).Value = 6.3  # T14: 6 -> 6.3
$ws.Cells.Item(14, 21).Value = 10.5  # U14: 9.75 -> 10.5
$ws.Cells.Item(14, 22).Value = 10  # V14: 9.75 -> 10
$ws.Cells.Item(14, 23).Value = 25  # W14: 23 -> 25
$ws.Cells.Item(14, 24).Value = 24  # X14: 23 -> 24
$ws.Cells.Item(14, 27).Value = 5.8  # AA14: 5.9 -> 5.8
$ws.Cells.Item(14, 31).Value = 7.2  # AE14: 7.5 -> 7.2
$ws.Cells.Item(14, 32).Value = 13.5  # AF14: 15 -> 13.5
$ws.Cells.Item(14, 33).Value = 11.25  # AG14: 12 -> 11.25
$ws.Cells.Item(14, 34).Value = 37  # AH14: 45 -> 37
$ws.Cells.Item(14, 35).Value = 32  # AI14: 35 -> 32
# Row 16
$ws.Cells.Item(16, 7).Value = 2.57  # G16: 3.15 -> 2.57
$ws.Cells.Item(16, 8).Value = 2.95  # H16: 3 -> 2.95
$ws.Cells.Item(16, 9).Value = 2.62  # I16: 2.18 -> 2.62
$ws.Cells.Item(16, 14).Value = 2.05  # N16: 2.02 -> 2.05
$ws.Cells.Item(16, 15).Value = 1.6  # O16: 1.62 -> 1.6
$ws.Cells.Item(16, 20).Value = 6.7  # T16: 7.6 -> 6.7
$ws.Cells.Item(16, 21).Value = 10.5  # U16: 13.5 -> 10.5
$ws.Cells.Item(16, 22).Value = 8.25  # V16: 9.25 -> 8.25
$ws.Cells.Item(16, 23).Value = 23  # W16: 32 -> 23
$ws.Cells.Item(16, 24).Value = 18  # X16: 23 -> 18
$ws.Cells.Item(16, 25).Value = 25  # Y16: 28 -> 25
$ws.Cells.Item(16, 26).Value = 7.8  # Z16: 8 -> 7.8
$ws.Cells.Item(16, 27).Value = 5.1  # AA16: 5.2 -> 5.1
$ws.Cells.Item(16, 28).Value = 11.5  # AB16: 11.75 -> 11.5
$ws.Cells.Item(16, 29).Value = 50  # AC16: 55 -> 50
$ws.Cells.Item(16, 31).Value = 6.3  # AE16: 5.8 -> 6.3
$ws.Cells.Item(16, 32).Value = 10.25  # AF16: 8.25 -> 10.25
$ws.Cells.Item(16, 33).Value = 8.5  # AG16: 7.6 -> 8.5
$ws.Cells.Item(16, 34).Value = 23  # AH16: 17 -> 23
$ws.Cells.Item(16, 35).Value = 19.5  # AI16: 15.5 -> 19.5
$ws.Cells.Item(16, 36).Value = 28  # AJ16: 25 -> 28
# Row 17
$ws.Cells.Item(17, 7).Value = 1.5  # G17: 1.47 -> 1.5
$ws.Cells.Item(17, 8).Value = 3.65  # H17: 3.8 -> 3.65
$ws.Cells.Item(17, 9).Value = 5.9  # I17: 6.1 -> 5.9
$ws.Cells.Item(17, 14).Value = 1.93  # N17: 1.91 -> 1.93
$ws.Cells.Item(17, 17).Value = 2.52  # Q17: 2.5 -> 2.52
$ws.Cells.Item(17, 20).Value = 5  # T17: 4.9 -> 5
$ws.Cells.Item(17, 21).Value = 5.4  # U17: 5.3 -> 5.4
$ws.Cells.Item(17, 22).Value = 7  # V17: 7.1 -> 7
$ws.Cells.Item(17, 23).Value = 8.25  # W17: 8 -> 8.25
$ws.Cells.Item(17, 26).Value = 8.5  # Z17: 9 -> 8.5
$ws.Cells.Item(17, 27).Value = 6.4  # AA17: 6.6 -> 6.4
$ws.Cells.Item(17, 31).Value = 11  # AE17: 11.75 -> 11
$ws.Cells.Item(17, 32).Value = 27  # AF17: 29 -> 27
$ws.Cells.Item(17, 33).Value = 16  # AG17: 16.5 -> 16
$ws.Cells.Item(17, 36).Value = 60  # AJ17: 55 -> 60
# Row 20
$ws.Cells.Item(20, 7).Value = 3.25  # G20: 3.75 -> 3.25
$ws.Cells.Item(20, 8).Value = 4  # H20: 4.25 -> 4
$ws.Cells.Item(20, 9).Value = 1.91  # I20: 1.72 -> 1.91
$ws.Cells.Item(20, 12).Value = 1.12  # L20: 1.11 -> 1.12
$ws.Cells.Item(20, 13).Value = 5.4  # M20: 5.6 -> 5.4
$ws.Cells.Item(20, 14).Value = 1.37  # N20: 1.36 -> 1.37
$ws.Cells.Item(20, 15).Value = 2.85  # O20: 2.9 -> 2.85
$ws.Cells.Item(20, 16).Value = 1.23  # P20: 1.21 -> 1.23
$ws.Cells.Item(20, 17).Value = 3.8  # Q20: 3.9 -> 3.8
$ws.Cells.Item(20, 18).Value = 1.38  # R20: 1.39 -> 1.38
$ws.Cells.Item(20, 19).Value = 2.82  # S20: 2.75 -> 2.82
$ws.Cells.Item(20, 20).Value = 19.5  # T20: 21 -> 19.5
$ws.Cells.Item(20, 21).Value = 25  # U20: 29 -> 25
$ws.Cells.Item(20, 22).Value = 12  # V20: 13.5 -> 12
$ws.Cells.Item(20, 23).Value = 50  # W20: 60 -> 50
$ws.Cells.Item(20, 24).Value = 23  # X20: 27 -> 23
$ws.Cells.Item(20, 25).Value = 21  # Y20: 24 -> 21
$ws.Cells.Item(20, 26).Value = 26  # Z20: 27 -> 26
$ws.Cells.Item(20, 27).Value = 9  # AA20: 9.5 -> 9
$ws.Cells.Item(20, 28).Value = 11  # AB20: 11.75 -> 11
$ws.Cells.Item(20, 29).Value = 28  # AC20: 30 -> 28
$ws.Cells.Item(20, 30).Value = 120  # AD20: 150 -> 120
$ws.Cells.Item(20, 32).Value = 13.5  # AF20: 12.5 -> 13.5
$ws.Cells.Item(20, 33).Value = 9  # AG20: 8.75 -> 9
$ws.Cells.Item(20, 34).Value = 19.5  # AH20: 17 -> 19.5
$ws.Cells.Item(20, 35).Value = 13  # AI20: 11.75 -> 13
$ws.Cells.Item(20, 36).Value = 16  # AJ20: 15.5 -> 16
# Row 22
$ws.Cells.Item(22, 13).Value = 4.5  # M22: 4.45 -> 4.5
$ws.Cells.Item(22, 20).Value = 11  # T22: 11.5 -> 11
$ws.Cells.Item(22, 21).Value = 12.5  # U22: 13 -> 12.5
$ws.Cells.Item(22, 25).Value = 18  # Y22: 17.5 -> 18
$ws.Cells.Item(22, 29).Value = 32  # AC22: 35 -> 32
$ws.Cells.Item(22, 31).Value = 15  # AE22: 14 -> 15
$ws.Cells.Item(22, 32).Value = 22  # AF22: 21 -> 22
$ws.Cells.Item(22, 35).Value = 24  # AI22: 25 -> 24
$ws.Cells.Item(22, 36).Value = 23  # AJ22: 25 -> 23
# Row 23
$ws.Cells.Item(23, 7).Value = 2.37  # G23: 2.55 -> 2.37
$ws.Cells.Item(23, 8).Value = 3.5  # H23: 3.55 -> 3.5
$ws.Cells.Item(23, 9).Value = 2.67  # I23: 2.45 -> 2.67
$ws.Cells.Item(23, 11).Value = 9  # K23: 9.25 -> 9
$ws.Cells.Item(23, 12).Value = 1.17  # L23: 1.16 -> 1.17
$ws.Cells.Item(23, 13).Value = 4.5  # M23: 4.55 -> 4.5
$ws.Cells.Item(23, 14).Value = 1.52  # N23: 1.5 -> 1.52
$ws.Cells.Item(23, 15).Value = 2.37  # O23: 2.4 -> 2.37
$ws.Cells.Item(23, 16).Value = 1.29  # P23: 1.28 -> 1.29
$ws.Cells.Item(23, 17).Value = 3.3  # Q23: 3.35 -> 3.3
$ws.Cells.Item(23, 19).Value = 2.62  # S23: 2.6 -> 2.62
$ws.Cells.Item(23, 20).Value = 12  # T23: 13.5 -> 12
$ws.Cells.Item(23, 21).Value = 15  # U23: 17 -> 15
$ws.Cells.Item(23, 22).Value = 9.25  # V23: 9.75 -> 9.25
$ws.Cells.Item(23, 23).Value = 27  # W23: 32 -> 27
$ws.Cells.Item(23, 24).Value = 17  # X23: 18 -> 17
$ws.Cells.Item(23, 25).Value = 19.5  # Y23: 20 -> 19.5
$ws.Cells.Item(23, 26).Value = 9  # Z23: 9.25 -> 9
$ws.Cells.Item(23, 27).Value = 7.4  # AA23: 7.5 -> 7.4
$ws.Cells.Item(23, 28).Value = 10.5  # AB23: 10.75 -> 10.5
$ws.Cells.Item(23, 30).Value = 150  # AD23: 175 -> 150
$ws.Cells.Item(23, 31).Value = 14  # AE23: 12.5 -> 14
$ws.Cells.Item(23, 32).Value = 18.5  # AF23: 15.5 -> 18.5
$ws.Cells.Item(23, 33).Value = 10  # AG23: 9.5 -> 10
$ws.Cells.Item(23, 34).Value = 35  # AH23: 28 -> 35
$ws.Cells.Item(23, 35).Value = 18.5  # AI23: 17.5 -> 18.5
# Row 24
$ws.Cells.Item(24, 8).Value = 4  # H24: 3.95 -> 4
$ws.Cells.Item(24, 9).Value = 2.42  # I24: 2.45 -> 2.42
$ws.Cells.Item(24, 13).Value = 6.5  # M24: 6.4 -> 6.5
$ws.Cells.Item(24, 14).Value = 1.27  # N24: 1.28 -> 1.27
$ws.Cells.Item(24, 15).Value = 3.45  # O24: 3.35 -> 3.45
$ws.Cells.Item(24, 17).Value = 4.15  # Q24: 4.2 -> 4.15
$ws.Cells.Item(24, 18).Value = 1.27  # R24: 1.28 -> 1.27
$ws.Cells.Item(24, 19).Value = 3.45  # S24: 3.4 -> 3.45
$ws.Cells.Item(24, 20).Value = 21  # T24: 20 -> 21
$ws.Cells.Item(24, 21).Value = 22  # U24: 21 -> 22
$ws.Cells.Item(24, 22).Value = 11  # V24: 10.75 -> 11
$ws.Cells.Item(24, 27).Value = 10.5  # AA24: 10 -> 10.5
$ws.Cells.Item(24, 29).Value = 21  # AC24: 22 -> 21
$ws.Cells.Item(24, 30).Value = 75  # AD24: 80 -> 75
$ws.Cells.Item(24, 31).Value = 22  # AE24: 21 -> 22
$ws.Cells.Item(24, 32).Value = 23  # AF24: 22 -> 23
$ws.Cells.Item(24, 35).Value = 16.5  # AI24: 17 -> 16.5
$ws.Cells.Item(24, 36).Value = 15.5  # AJ24: 16 -> 15.5
# Row 25
$ws.Cells.Item(25, 7).Value = 2.3  # G25: 2.5 -> 2.3
$ws.Cells.Item(25, 9).Value = 3.4  # I25: 3 -> 3.4
$ws.Cells.Item(25, 12).Value = 1.36  # L25: 1.4 -> 1.36
$ws.Cells.Item(25, 13).Value = 3  # M25: 2.75 -> 3
$ws.Cells.Item(25, 16).Value = 1.5  # P25: 1.44 -> 1.5
$ws.Cells.Item(25, 17).Value = 2.5  # Q25: 2.63 -> 2.5
$ws.Cells.Item(25, 18).Value = 1.95  # R25: 1.91 -> 1.95
$ws.Cells.Item(25, 19).Value = 1.8  # S25: 1.91 -> 1.8
$ws.Cells.Item(25, 20).Value = 7  # T25: 7.5 -> 7
$ws.Cells.Item(25, 21).Value = 10  # U25: 12 -> 10
$ws.Cells.Item(25, 22).Value = 9.5  # V25: 10 -> 9.5
$ws.Cells.Item(25, 23).Value = 21  # W25: 23 -> 21
$ws.Cells.Item(25, 26).Value = 7.5  # Z25: 8 -> 7.5
$ws.Cells.Item(25, 30).Value = 351  # AD25: 301 -> 351
$ws.Cells.Item(25, 31).Value = 9  # AE25: 8.5 -> 9
$ws.Cells.Item(25, 32).Value = 15  # AF25: 13 -> 15
$ws.Cells.Item(25, 33).Value = 12  # AG25: 11 -> 12
$ws.Cells.Item(25, 34).Value = 34  # AH25: 29 -> 34
$ws.Cells.Item(25, 35).Value = 29  # AI25: 26 -> 29
$ws.Cells.Item(25, 36).Value = 41  # AJ25: 34 -> 41
# Row 30
$ws.Cells.Item(30, 7).Value = 1.5  # G30: 1.53 -> 1.5
